# Auto-generated edit script
# Updates currentAveragePrice / profit columns (H-N) across 8 crafting-job
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) following a scheduled
# market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1018.04
$ws.Range("I80").Value = 621.1
$ws.Range("J80").Value = 1282.6666
$ws.Range("K80").Value = 1863.3
$ws.Range("L80").Value = 3847.9998
$ws.Range("M80").Value = -865.3000000000002
$ws.Range("N80").Value = -5843.9998
$ws.Range("H83").Value = 1018.04
$ws.Range("I83").Value = 621.1
$ws.Range("J83").Value = 1282.6666
$ws.Range("K83").Value = 5589.900000000001
$ws.Range("L83").Value = 11543.9994
$ws.Range("M83").Value = -597.9000000000005
$ws.Range("N83").Value = -21527.9994
$ws.Range("H86").Value = 2494.5
$ws.Range("I86").Value = 1174.5714
$ws.Range("J86").Value = 3814.4285
$ws.Range("K86").Value = 1174.5714
$ws.Range("L86").Value = 3814.4285
$ws.Range("M86").Value = -51.57140000000004
$ws.Range("N86").Value = -6060.4285
$ws.Range("H89").Value = 2494.5
$ws.Range("I89").Value = 1174.5714
$ws.Range("J89").Value = 3814.4285
$ws.Range("K89").Value = 5872.857
$ws.Range("L89").Value = 19072.1425
$ws.Range("M89").Value = -256.857
$ws.Range("N89").Value = -30304.1425
$ws.Range("H111").Value = 8730.807000000001
$ws.Range("I111").Value = 11589.789
$ws.Range("J111").Value = 4204.0835
$ws.Range("K111").Value = 34769.367
$ws.Range("L111").Value = 12612.2505
$ws.Range("M111").Value = -31702.367
$ws.Range("N111").Value = -18746.2505
$ws.Range("H125").Value = 10000
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H132").Value = 7037.82
$ws.Range("I132").Value = 2432.3044
$ws.Range("J132").Value = 60001.25
$ws.Range("K132").Value = 7296.9132
$ws.Range("L132").Value = 180003.75
$ws.Range("M132").Value = -4766.9132
$ws.Range("N132").Value = -185063.75
$ws.Range("H137").Value = 3185
$ws.Range("I137").Value = 3511.9048
$ws.Range("K137").Value = 10535.7144
$ws.Range("M137").Value = -7985.714399999999
$ws.Range("H138").Value = 214507.02
$ws.Range("I138").Value = 1063.8572
$ws.Range("K138").Value = 3191.5716
$ws.Range("M138").Value = 1948.4284
$ws.Range("H141").Value = 3424.476
$ws.Range("I141").Value = 2526.3125
$ws.Range("K141").Value = 7578.9375
$ws.Range("M141").Value = -2398.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1346.3636
$ws.Range("I2").Value = 1388.5
$ws.Range("K2").Value = 1388.5
$ws.Range("M2").Value = -1275.5
$ws.Range("H61").Value = 6214.3784
$ws.Range("I61").Value = 6321.5586
$ws.Range("K61").Value = 6321.5586
$ws.Range("M61").Value = -6109.5586
$ws.Range("H63").Value = 2643.3333
$ws.Range("I63").Value = 2481.6667
$ws.Range("J63").Value = 2966.6667
$ws.Range("K63").Value = 2481.6667
$ws.Range("L63").Value = 2966.6667
$ws.Range("M63").Value = -1795.6667
$ws.Range("N63").Value = -4338.6667
$ws.Range("H66").Value = 2643.3333
$ws.Range("I66").Value = 2481.6667
$ws.Range("J66").Value = 2966.6667
$ws.Range("K66").Value = 12408.3335
$ws.Range("L66").Value = 14833.3335
$ws.Range("M66").Value = -8976.333500000001
$ws.Range("N66").Value = -21697.3335
$ws.Range("H102").Value = 16668555
$ws.Range("I102").Value = 2061.5557
$ws.Range("K102").Value = 2061.5557
$ws.Range("M102").Value = -439.5556999999999
$ws.Range("H109").Value = 85000
$ws.Range("J109").Value = 85000
$ws.Range("L109").Value = 85000
$ws.Range("N109").Value = -87774
$ws.Range("H116").Value = 1346.3636
$ws.Range("I116").Value = 1388.5
$ws.Range("K116").Value = 1388.5
$ws.Range("M116").Value = 905.5
$ws.Range("H132").Value = 2800.7188
$ws.Range("I132").Value = 2704.8035
$ws.Range("J132").Value = 3472.125
$ws.Range("K132").Value = 8114.4105
$ws.Range("L132").Value = 10416.375
$ws.Range("M132").Value = -5584.4105
$ws.Range("N132").Value = -15476.375
$ws.Range("H136").Value = 6214.3784
$ws.Range("I136").Value = 6321.5586
$ws.Range("K136").Value = 18964.6758
$ws.Range("M136").Value = -16414.6758

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1346.3636
$ws.Range("I3").Value = 1388.5
$ws.Range("K3").Value = 1388.5
$ws.Range("M3").Value = -1274.5
$ws.Range("H107").Value = 2169.5386
$ws.Range("I107").Value = 1701
$ws.Range("J107").Value = 2571.1428
$ws.Range("K107").Value = 1701
$ws.Range("L107").Value = 2571.1428
$ws.Range("M107").Value = 219
$ws.Range("N107").Value = -6411.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 19980
$ws.Range("I39").Value = 4950
$ws.Range("J39").Value = 30000
$ws.Range("K39").Value = 4950
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = -4559
$ws.Range("N39").Value = -30782
$ws.Range("H49").Value = 19980
$ws.Range("I49").Value = 4950
$ws.Range("J49").Value = 30000
$ws.Range("K49").Value = 4950
$ws.Range("L49").Value = 30000
$ws.Range("M49").Value = -4768
$ws.Range("N49").Value = -30364
$ws.Range("H58").Value = 2296.8
$ws.Range("I58").Value = 2495.5715
$ws.Range("K58").Value = 2495.5715
$ws.Range("M58").Value = -2292.5715
$ws.Range("H134").Value = 2866.8157
$ws.Range("I134").Value = 1386.7407
$ws.Range("K134").Value = 4160.2221
$ws.Range("M134").Value = -1625.2221
$ws.Range("H136").Value = 2296.8
$ws.Range("I136").Value = 2495.5715
$ws.Range("K136").Value = 7486.7145
$ws.Range("M136").Value = -4936.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 272.5
$ws.Range("I19").Value = 272.5
$ws.Range("K19").Value = 817.5
$ws.Range("M19").Value = -643.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5379617
$ws.Range("I132").Value = 6413399.5
$ws.Range("J132").Value = 3946.5
$ws.Range("K132").Value = 19240198.5
$ws.Range("L132").Value = 11839.5
$ws.Range("M132").Value = -19237668.5
$ws.Range("N132").Value = -16899.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5925.2085
$ws.Range("J40").Value = 5872
$ws.Range("L40").Value = 5872
$ws.Range("N40").Value = -6144
$ws.Range("H132").Value = 3065.8774
$ws.Range("I132").Value = 3059.7026
$ws.Range("J132").Value = 3084.9167
$ws.Range("K132").Value = 9179.1078
$ws.Range("L132").Value = 9254.750100000001
$ws.Range("M132").Value = -6649.1078
$ws.Range("N132").Value = -14314.7501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 12476.579
$ws.Range("I136").Value = 15862.357
$ws.Range("K136").Value = 47587.071
$ws.Range("M136").Value = -45037.071
